# "Final version of dataset" - update correlation-matrix inputs on the
# "cor" sheet (rows 12-26) and leave the active selection on R19, matching
# where the author's cursor ended up when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cor")

# Row 12
$ws.Range("A12").Value = 0.2
$ws.Range("B12").Value = -0.2

# Row 13
$ws.Range("A13").Value = 0.2

# Row 14
$ws.Range("A14").Value = 0.2
$ws.Range("B14").Value = -0.2

# Row 15
$ws.Range("A15").Value = 0
$ws.Range("B15").Value = -0.2
$ws.Range("C15").Value = 0.2
$ws.Range("D15").Value = 0.2
$ws.Range("E15").Value = 0.1
$ws.Range("F15").Value = 0.1

# Row 16
$ws.Range("A16").Value = 0
$ws.Range("B16").Value = -0.2
$ws.Range("C16").Value = 0.2
$ws.Range("D16").Value = 0.2
$ws.Range("E16").Value = 0.1
$ws.Range("F16").Value = 0.1

# Row 17
$ws.Range("A17").Value = 0
$ws.Range("B17").Value = -0.2
$ws.Range("C17").Value = 0.2
$ws.Range("D17").Value = 0.2
$ws.Range("E17").Value = 0.1
$ws.Range("F17").Value = 0.1

# Row 18
$ws.Range("A18").Value = -0.2
$ws.Range("B18").Value = -0.2
$ws.Range("C18").Value = 0.2
$ws.Range("D18").Value = 0.2
$ws.Range("E18").Value = 0.1
$ws.Range("F18").Value = 0.1

# Row 19
$ws.Range("A19").Value = -0.2
$ws.Range("B19").Value = -0.2
$ws.Range("C19").Value = 0.2
$ws.Range("D19").Value = 0.2
$ws.Range("E19").Value = 0.1
$ws.Range("F19").Value = 0.1
$ws.Range("N19").Value = 0.1
$ws.Range("O19").Value = 0.1
$ws.Range("P19").Value = 0.1
$ws.Range("Q19").Value = 0.1

# Row 20
$ws.Range("A20").Value = -0.2
$ws.Range("B20").Value = -0.2
$ws.Range("C20").Value = 0.2
$ws.Range("D20").Value = 0.2
$ws.Range("E20").Value = 0.1
$ws.Range("F20").Value = 0.1
$ws.Range("O20").Value = 0.1
$ws.Range("P20").Value = 0.1
$ws.Range("Q20").Value = 0.1

# Row 21
$ws.Range("A21").Value = -0.2
$ws.Range("B21").Value = -0.2
$ws.Range("C21").Value = 0.2
$ws.Range("D21").Value = 0.2
$ws.Range("E21").Value = 0.1
$ws.Range("F21").Value = 0.1
$ws.Range("O21").Value = 0.1
$ws.Range("P21").Value = 0.1
$ws.Range("Q21").Value = 0.1

# Row 22
$ws.Range("A22").Value = -0.2
$ws.Range("B22").Value = -0.2
$ws.Range("C22").Value = 0.2
$ws.Range("D22").Value = 0.2
$ws.Range("E22").Value = 0.1
$ws.Range("F22").Value = 0.1
$ws.Range("O22").Value = 0.1
$ws.Range("P22").Value = 0.1
$ws.Range("Q22").Value = 0.1

# Row 23
$ws.Range("A23").Value = -0.1
$ws.Range("B23").Value = -0.2
$ws.Range("C23").Value = 0.2
$ws.Range("D23").Value = 0.2
$ws.Range("E23").Value = 0.1
$ws.Range("F23").Value = 0.1

# Row 24
$ws.Range("A24").Value = -0.1
$ws.Range("B24").Value = -0.2
$ws.Range("C24").Value = 0.2
$ws.Range("D24").Value = 0.2
$ws.Range("E24").Value = 0.1
$ws.Range("F24").Value = 0.1

# Row 25
$ws.Range("A25").Value = -0.1
$ws.Range("B25").Value = -0.2
$ws.Range("C25").Value = 0.2
$ws.Range("D25").Value = 0.2
$ws.Range("E25").Value = 0.1
$ws.Range("F25").Value = 0.1

# Row 26
$ws.Range("S26").Value = 0.5
$ws.Range("T26").Value = 0.5
$ws.Range("U26").Value = 0.5

# Leave the cursor where the author left it when the workbook was saved.
$ws.Range("R19").Select() | Out-Null
